$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up the common-name column (C): drop the trailing spaces / plural
# forms and replace with the tidy singular common names.
$ws.Range("C2").Value = "Grape"
$ws.Range("C3").Value = "Blueberry"
$ws.Range("C4").Value = "Apple"
$ws.Range("C5").Value = "Cherry"
$ws.Range("C6").Value = "Almond"
$ws.Range("C7").Value = "Walnut"

# These common-name cells were italicised - remove that formatting so the
# whole column reads consistently (C2 and C4 were already plain).
$ws.Range("C3").Font.Italic = $false
$ws.Range("C5").Font.Italic = $false
$ws.Range("C6").Font.Italic = $false
$ws.Range("C7").Font.Italic = $false

# Move the selection to C2, and park the window where it last was left.
$ws.Range("C2").Select()
